$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ScopePropertiesApi"
$ws.Range("B2").Value = "ScopePropertiesController"
$ws.Range("C2").Value = "Get"
$ws.Range("E2").Value = "A"
$ws.Range("F2").Value = "Claims"
$ws.Range("G2").Value = "claim*name=jack&claim*role=admin"

$ws.Range("A3").Value = "ScopePropertiesApi"
$ws.Range("B3").Value = "ScopePropertiesController"
$ws.Range("C3").Value = "Get"
$ws.Range("E3").Value = "A"
$ws.Range("F3").Value = "Headers"
$ws.Range("G3").Value = "header*hdr1=ABC&header*hdr2=DEF"

$ws.Range("A4").Value = "ScopePropertiesApi"
$ws.Range("B4").Value = "ScopePropertiesController"
$ws.Range("C4").Value = "Get"
$ws.Range("E4").Value = "A"
$ws.Range("F4").Value = "Expected"
$ws.Range("G4").Value = "{`"User`":`"jack`",`"name`":`"jack`",`"role`":`"admin`",`"Host`":`"localhost`",`"hdr1`":`"ABC`",`"hdr2`":`"DEF`",`"X-HostPath`":`"localhost`"}"

$ws.Range("A5").Value = "ScopePropertiesApi"
$ws.Range("B5").Value = "ScopePropertiesController"
$ws.Range("C5").Value = "Get"
$ws.Range("E5").Value = "B"
$ws.Range("F5").Value = "Claims"
$ws.Range("G5").Value = "claim*role=user&claim*group=456"

$ws.Range("A6").Value = "ScopePropertiesApi"
$ws.Range("B6").Value = "ScopePropertiesController"
$ws.Range("C6").Value = "Get"
$ws.Range("E6").Value = "B"
$ws.Range("F6").Value = "Headers"
$ws.Range("G6").Value = "header*hdr1=123&header*X-User=jill"

$ws.Range("A7").Value = "ScopePropertiesApi"
$ws.Range("B7").Value = "ScopePropertiesController"
$ws.Range("C7").Value = "Get"
$ws.Range("E7").Value = "B"
$ws.Range("F7").Value = "Expected"
$ws.Range("G7").Value = "{`"User`":`"jill`",`"role`":`"user`",`"group`":`"456`",`"Host`":`"localhost`",`"hdr1`":`"123`",`"X-User`":`"jill`"}"

$ws.Range("A8").Value = "ScopePropertiesApi"
$ws.Range("B8").Value = "ScopePropertiesController"
$ws.Range("C8").Value = "Get"
$ws.Range("E8").Value = "C"
$ws.Range("F8").Value = "Claims"
$ws.Range("G8").Value = "X-User=bob"

$ws.Range("A9").Value = "ScopePropertiesApi"
$ws.Range("B9").Value = "ScopePropertiesController"
$ws.Range("C9").Value = "Get"
$ws.Range("E9").Value = "C"
$ws.Range("F9").Value = "Headers"
$ws.Range("G9").Value = "header*hdr1=123&header*X-User=jill"

$ws.Range("A10").Value = "ScopePropertiesApi"
$ws.Range("B10").Value = "ScopePropertiesController"
$ws.Range("C10").Value = "Get"
$ws.Range("E10").Value = "C"
$ws.Range("F10").Value = "Expected"
$ws.Range("G10").Value = "{`"User`":`"bob`"}"

$ws.Range("G10").Select() | Out-Null
